$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: only the row height changes (data/styles unchanged) ---
$ws.Rows.Item(9).RowHeight = 28.5

# --- Row 11: new problem entry "Valid Sudoku" ---
# (filled in before row 10's edit so the new "Valid Sudoku" shared string
#  is allocated index 46, matching the diff's string-table order)
$ws.Rows.Item(11).RowHeight = 28.5
$ws.Range("A11").Value = 36
$ws.Range("B11").Value = "Valid Sudoku"
$ws.Range("C11").Value = "Array"
$ws.Range("D11").Value = "Array, Hash table, matrix"
$ws.Range("F11").Value = "Medium"
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = "✅"
$ws.Range("I11").Value = "Copied from Solution pane and didn't understand"

# Match formatting to the sibling rows that already use these same styles
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)   # xlPasteFormats - red fill (s=5)

$ws.Range("H10").Copy()
$ws.Range("H11").PasteSpecial(-4122)   # xlPasteFormats - red Calibri font (s=2)

$ws.Range("D10").Copy()
$ws.Range("D11").PasteSpecial(-4122)   # xlPasteFormats - wrap text (s=8)

$ws.Range("I10").Copy()
$ws.Range("I11").PasteSpecial(-4122)   # xlPasteFormats - wrap text (s=8)

# --- Row 10: height change + D10 tag text change ---
$ws.Rows.Item(10).RowHeight = 33
$ws.Range("D10").Value = "Array, Binary search, matrix"

# --- Sheet view: selection moves to B15 (also clears the stale
#     topLeftCell="A5" scroll-position left over from the previous edit) ---
$ws.Range("B15").Select() | Out-Null

Write-Host "done"
